$wb = $excel.ActiveWorkbook

$mdUrl     = "https://github.com/OpenLocalizationTest/oltest/blob/ae112ff2954990959451a2469ba1a280987d7f4f/e2e/aa6ca3cb-c3ce-4575-8446-36a6625c758e.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ae112ff2954990959451a2469ba1a280987d7f4f/.localization-config"

# The Overview sheet mirrors the per-language "Status" text via the same shared
# string, so it must be updated too or the old "Not yet handed off" string will
# stick around in the shared string table.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff failed"
$wsOverview.Range("C2").Value = "Handoff failed"

$hyperColor = 15570276   # BGR for RGB FF6495ED (matches workbook's custom HyperLink style)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # --- Update row 2 (handoff now failed / ignored, no more "Latest Handoff File") ---
    $ws.Range("B2").Value = "Handoff failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    # --- Rebuild hyperlinks, dropping the old "Latest Handoff File" (C2) link ---
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "aa6ca3cb-c3ce-4575-8446-36a6625c758e.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config")

    # Restore the original hyperlink-cell look (the Add() call above applies Excel's
    # default "Hyperlink" look-and-feel, so bring the font back in line with the
    # workbook's own custom HyperLink style).
    $ws.Range("A2").Font.Color = $hyperColor
    $ws.Range("A2").Font.Underline = 2
    $ws.Range("A3").Font.Color = $hyperColor
    $ws.Range("A3").Font.Underline = 2
}
